$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.586.66"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.65%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.820.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.54%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("E5").Value = "  +0.15%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "305.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.53%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4666"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3591"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.97%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07126"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.46%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9020"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.71%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07800"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.71%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.67%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.839.59"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.39%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.251"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.39%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.327"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.06%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.22"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.51%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.010"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.14%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008553"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.008"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.09%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.648.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.66%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.007"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.54%  "

$ws.Range("E23").Value = "  +0.64%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.944"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.82%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.45%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.967"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.47%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "113.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.28%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.798"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.92%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08803"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.67%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.145"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.79%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.765"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.93%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7280"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.73%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.437"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.20%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.122"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.86%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.074"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.23%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01925"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.80%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.914"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.67%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05105"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.01%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.835"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.86%  "

$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5047"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.85%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1494"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.98%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.979"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.41%  "

$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.009"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.17%  "

$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4653"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.71%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.956"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.16%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "98.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.557"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.68%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06001"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.65%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "63.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.61%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "35.66"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.48%  "
